$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column K corresponds to "municipio-nombre" metadata. With the new curated
# dimensions, it is re-classified from a "measure" into a "dimension"
# (refArea), using a URI-Municipio identifier pattern instead of xsd:int.
$ws.Range("K2").Value = "sdmx-dimension:refArea"
$ws.Range("K3").Value = "dim"
$ws.Range("K4").Value = "URI-Municipio"
